$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    "A1" = 2504.3641874528048
    "B1" = 1675.4872758024103
    "C1" = 1659.1694206955858
    "A2" = 2263.2649304931097
    "B2" = 1537.5578724368795
    "C2" = 1393.2740556789715
    "A3" = 2560.3680387196541
    "B3" = 1740.6137032078921
    "C3" = 1579.1404065201828
    "A4" = 2495.9933054598469
    "B4" = 1913.3276998137765
    "C4" = 1951.620215286079
    "A5" = 2528.8045935782193
    "B5" = 1752.3685036322445
    "C5" = 1755.617454012624
    "A6" = 2480.0542736911661
    "B6" = 1850.23187819243
    "C6" = 1930.686298371365
    "A7" = 2397.2088804141163
    "B7" = 1848.4313192171414
    "C7" = 1675.6626888020246
    "A8" = 2467.9969994687162
    "B8" = 1930.0304952068079
    "C8" = 1785.2212237624467
    "A9" = 2643.8787844459366
    "B9" = 1957.2149355131842
    "C9" = 1663.4754274722636
    "A10" = 2394.9100079198947
    "B10" = 1513.0137755996332
    "C10" = 1430.9226292086887
    "A11" = 2163.0420354858652
    "B11" = 1575.4241850871933
    "C11" = 1407.2608605798839
    "A12" = 2770.6265071879297
    "B12" = 2201.7324782617106
    "C12" = 1896.5290666105197
    "A13" = 2528.2387483210841
    "B13" = 1947.5471764946674
    "C13" = 1749.6320391353865
    "A14" = 2609.7077576152847
    "B14" = 2018.747131381228
    "C14" = 1757.4978996681295
    "A15" = 2503.2248039583947
    "B15" = 2051.8787234510228
    "C15" = 1832.9544835080455
    "A16" = 2592.376243186462
    "B16" = 1799.9947305053117
    "C16" = 1571.4178997318688
    "A17" = 2356.8601378921621
    "B17" = 1814.183856138367
    "C17" = 1709.7084558778279
    "A18" = 2616.1384098191234
    "B18" = 2161.0004518584929
    "C18" = 2088.9723595506885
    "A19" = 2011.2295767214825
    "B19" = 2033.9315696339561
    "C19" = 1981.0685119676493
    "A20" = 2576.2845538913493
    "B20" = 1981.4044547233398
    "C20" = 1869.9530479197674
    "A21" = 2713.0801992948718
    "B21" = 2020.856815124363
    "C21" = 1980.6745582119877
    "A22" = 2569.5114393257518
    "B22" = 1982.44279635917
    "C22" = 1731.9157297969646
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
